# Update answer of chap8 and chap9
# Target: ppt/slides/slide1.xml (the "chapter8-image" connectors diagram)
#
# The PowerPoint COM object model expresses shape position/size in points
# (stored internally as 32-bit floats), while the underlying OOXML stores
# EMU (1 pt = 12700 EMU) and truncates the float32-roundtripped point value
# back to an integer EMU on save. A plain `emu / 12700.0` division can land
# a few EMU below the intended integer once it has been narrowed to a
# float32, so PtForEmu nudges the point value upward (in tiny steps) until
# the float32 round-trip truncates back to exactly the target EMU.

$EMU_PER_PT = 12700.0

function PtForEmu([double]$targetEmu) {
    $base = $targetEmu / $EMU_PER_PT
    $pt = $base
    for ($i = 0; $i -lt 4096; $i++) {
        $f32 = [single]$pt
        $back = [double]$f32 * $EMU_PER_PT
        $flo = [math]::Floor($back + 0.0000001)
        if ([long]$flo -eq [long]$targetEmu) {
            return $pt
        }
        $pt = $pt + 0.0000002384185791015625
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Shrink the connector "直接连接符 24" (currently 999074 x 1592788 EMU)
#    down to 299430 x 477371 EMU, keeping its existing top-left offset.
# ---------------------------------------------------------------------
$shp24 = $s.Shapes.Item(6)
$shp24.Width = PtForEmu 299430
$shp24.Height = PtForEmu 477371

# ---------------------------------------------------------------------
# 2) Reposition/resize connector "直接连接符 52"
#    off: (3462619,1398494) -> (2951629,903620)
#    ext: (1055593,1022301) -> (1566584,1517176)
# ---------------------------------------------------------------------
$shp52 = $s.Shapes.Item(10)
$shp52.Left = PtForEmu 2951629
$shp52.Top = PtForEmu 903620
$shp52.Width = PtForEmu 1566584
$shp52.Height = PtForEmu 1517176

# ---------------------------------------------------------------------
# 3) Add a new straight connector "直接连接符 13" after "直接连接符 59".
#    Duplicate an existing unflipped, unstyled-override connector so the
#    clone inherits the same <p:style> (lnRef/fillRef/effectRef/fontRef)
#    and <a:cxnSpLocks/> structure used throughout this diagram, then
#    move/resize/rename it to match the target geometry.
# ---------------------------------------------------------------------
$template = $s.Shapes.Item(4)
$newRange = $template.Duplicate()
$newShape = $newRange.Item(1)

$newShape.Name = "直接连接符 13"
$newShape.Left = PtForEmu 6272657
$newShape.Top = PtForEmu 4255994
$newShape.Width = PtForEmu 1882984
$newShape.Height = PtForEmu 1385047

Write-Host "Done: shapes =" $s.Shapes.Count
